{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph separating it from the requirements line)\n// that followed the \"LOM3089: ... (Requisito fraco)\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph that holds the requirement line; the footer block\n// to delete is the three paragraphs that immediately follow it.\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Requisito fraco\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const toDelete = [];\n  for (let i = anchorIndex + 1; i < paragraphs.items.length; i++) {\n    const text = paragraphs.items[i].text;\n    if (\n      text === \"\" ||\n      text.indexOf(\"Ver no Jupiter\") !== -1 ||\n      text.indexOf(\"Powered by Jekyll\") !== -1\n    ) {\n      toDelete.push(paragraphs.items[i]);\n      if (text.indexOf(\"Powered by Jekyll\") !== -1) {\n        break;\n      }\n    } else {\n      break;\n    }\n  }\n  for (const para of toDelete) {\n    para.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph separating it from the requirements line)\n# that followed the \"LOM3089: ... (Requisito fraco)\" paragraph.\n$d = $word.ActiveDocument\n\n# Locate the paragraph holding the requirement line; the footer block to\n# remove is the run of paragraphs immediately following it: a blank\n# paragraph, the \"Ver no Jupiter...\" line, and the copyright line.\n$anchor = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Requisito fraco*\") {\n        $anchor = $i\n        break\n    }\n}\n\nif ($anchor -ne -1) {\n    $toDelete = @()\n    $i = $anchor + 1\n    while ($i -le $d.Paragraphs.Count) {\n        $t = $d.Paragraphs.Item($i).Range.Text\n        $trimmed = $t.Trim()\n        if ($trimmed -eq \"\" -or $t -like \"*Ver no Jupiter*\" -or $t -like \"*Powered by Jekyll*\") {\n            $toDelete += $i\n            if ($t -like \"*Powered by Jekyll*\") {\n                break\n            }\n            $i++\n        } else {\n            break\n        }\n    }\n\n    # Delete from the highest index down so earlier indices stay valid.\n    for ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n        $d.Paragraphs.Item($toDelete[$j]).Range.Delete()\n    }\n}\n"}
